# Update gh-pages output data for both the "展览" (exhibition) sheet and the
# "全部类型" (all-types) sheet, which carry identical tables.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Simple numeric refreshes on existing rows -----------------------
    $ws.Cells.Item(2, 6).Value = 1084     # F2  1081 -> 1084
    $ws.Cells.Item(5, 6).Value = 4641     # F5  4638 -> 4641
    $ws.Cells.Item(7, 6).Value = 390      # F7  389  -> 390
    $ws.Cells.Item(8, 6).Value = 1380     # F8  1378 -> 1380
    $ws.Cells.Item(11, 6).Value = 1099    # F11 1090 -> 1099
    $ws.Cells.Item(13, 6).Value = 597     # F13 588  -> 597
    $ws.Cells.Item(13, 7).Value = 58.5    # G13 60   -> 58.5
    $ws.Cells.Item(15, 6).Value = 24      # F15 20   -> 24

    # --- Make room for the new event row ---------------------------------
    # Before: row16=AP嘉年华(old)->CM01, row17=鹰潭. After: row16=new AP event
    # (small ticket), row17=CM01 (previously row16's content), row18=鹰潭
    # (previously row17's content). Insert a fresh blank row at the bottom
    # of the table (row 18) so rows 16/17 keep their identity while row 17's
    # event details are swapped out and the displaced data lands in new row 18.
    $ws.Rows.Item(18).Insert()

    # Give the new row's index cell (column A) the same bordered/centered
    # style used by every other row-index cell in the table.
    $ws.Cells.Item(17, 1).Copy()
    $ws.Cells.Item(18, 1).PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    # --- Row 18: gains the event data that used to live in row 17 --------
    $ws.Cells.Item(18, 1).Value = 17
    # Dotted dates like "2024.03.30" get auto-parsed into real Excel dates
    # unless the cell is forced to Text first; clear the format right back
    # afterwards so the cell ends up with no special style, same as source.
    $ws.Cells.Item(18, 2).NumberFormat = "@"
    $ws.Cells.Item(18, 2).Value = "2024.03.30"
    $ws.Cells.Item(18, 2).ClearFormats()
    $ws.Cells.Item(18, 3).Value = "鹰潭·原×铁×崩only"
    $ws.Cells.Item(18, 4).Value = "南站路24号 回禾酒店(鹰潭火车站店)"
    $ws.Cells.Item(18, 5).Value = "2024.03.30 10:00-03.30 17:00"
    $ws.Cells.Item(18, 6).Value = 23
    $ws.Cells.Item(18, 7).Value = 60
    $ws.Cells.Item(18, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81097"
    $ws.Cells.Item(18, 9).Value = "//i2.hdslb.com/bfs/openplatform/202401/q0AZaXAk1705646244207.jpeg"

    # --- Row 17: gains the event data that used to live in row 16 --------
    # (A17/B17 are left exactly as they already were: 16 / 2024.03.30)
    $ws.Cells.Item(17, 3).Value = "南昌·CM01动漫游戏博览会"
    $ws.Cells.Item(17, 4).Value = "怀玉山大道1315号 南昌绿地国际博览中心"
    $ws.Cells.Item(17, 5).Value = "2024.03.30 10:00-03.31 17:00"
    $ws.Cells.Item(17, 6).Value = 269
    $ws.Cells.Item(17, 7).Value = 55
    $ws.Cells.Item(17, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81691"
    $ws.Cells.Item(17, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/IYLaH7AS1706866218597.png"

    # --- Row 16: replaced with the brand-new event -------------------------
    # (A16 is left exactly as it already was: 15)
    $ws.Cells.Item(16, 2).NumberFormat = "@"
    $ws.Cells.Item(16, 2).Value = "2024.03.24"
    $ws.Cells.Item(16, 2).ClearFormats()
    $ws.Cells.Item(16, 3).Value = "南昌·AP动漫游戏  嘉年华内场票-小N&子音"
    $ws.Cells.Item(16, 4).Value = "八一桥街道青山南路118号 蓝海会展中心"
    $ws.Cells.Item(16, 5).Value = "2024.03.24 09:00-03.24 17:00"
    $ws.Cells.Item(16, 6).Value = 4
    $ws.Cells.Item(16, 7).Value = 218
    $ws.Cells.Item(16, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81973"
    $ws.Cells.Item(16, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/zbG5HICL1708504962467.jpeg"
}
